$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamps for A2:A97 (shifted from 2025-07-xx to 2025-07-xx + 12 days,
# i.e. a new day's worth of Entsoe solar generation forecast data).
$aValues = @(
    45875.01041666666, 45875.02083333334, 45875.03125,       45875.04166666666,
    45875.05208333334, 45875.0625,        45875.07291666666, 45875.08333333334,
    45875.09375,        45875.10416666666, 45875.11458333334, 45875.125,
    45875.13541666666, 45875.14583333334, 45875.15625,       45875.16666666666,
    45875.17708333334, 45875.1875,        45875.19791666666, 45875.20833333334,
    45875.21875,        45875.22916666666, 45875.23958333334, 45875.25,
    45875.26041666666, 45875.27083333334, 45875.28125,       45875.29166666666,
    45875.30208333334, 45875.3125,        45875.32291666666, 45875.33333333334,
    45875.34375,        45875.35416666666, 45875.36458333334, 45875.375,
    45875.38541666666, 45875.39583333334, 45875.40625,       45875.41666666666,
    45875.42708333334, 45875.4375,        45875.44791666666, 45875.45833333334,
    45875.46875,        45875.47916666666, 45875.48958333334, 45875.5,
    45875.51041666666, 45875.52083333334, 45875.53125,       45875.54166666666,
    45875.55208333334, 45875.5625,        45875.57291666666, 45875.58333333334,
    45875.59375,        45875.60416666666, 45875.61458333334, 45875.625,
    45875.63541666666, 45875.64583333334, 45875.65625,       45875.66666666666,
    45875.67708333334, 45875.6875,        45875.69791666666, 45875.70833333334,
    45875.71875,        45875.72916666666, 45875.73958333334, 45875.75,
    45875.76041666666, 45875.77083333334, 45875.78125,       45875.79166666666,
    45875.80208333334, 45875.8125,        45875.82291666666, 45875.83333333334,
    45875.84375,        45875.85416666666, 45875.86458333334, 45875.875,
    45875.88541666666, 45875.89583333334, 45875.90625,       45875.91666666666,
    45875.92708333334, 45875.9375,        45875.94791666666, 45875.95833333334,
    45875.96875,        45875.97916666666, 45875.98958333334, 45876
)

# New notified-production values for B18:B89 (all other B values are unchanged).
$bValues = @(
    1, 1, 1, 1,
    19, 20, 24, 32,
    179, 198, 218, 244,
    603, 636, 671, 710,
    1113, 1144, 1179, 1218,
    1471, 1500, 1524, 1548,
    1660, 1666, 1677, 1685,
    1733, 1734, 1733, 1730,
    1716, 1712, 1710, 1708,
    1647, 1638, 1625, 1607,
    1478, 1454, 1433, 1404,
    1218, 1190, 1161, 1135,
    894, 863, 839, 814,
    467, 443, 417, 396,
    157, 138, 124, 111,
    19, 13, 10, 10,
    8, 8, 8, 8,
    7, 7, 7, 7
)

$aRowCount = $aValues.Length
$aArr = New-Object 'object[,]' $aRowCount,1
for ($i = 0; $i -lt $aRowCount; $i++) {
    $aArr[$i,0] = $aValues[$i]
}
$ws.Range("A2:A97").Value2 = $aArr

$bRowCount = $bValues.Length
$bArr = New-Object 'object[,]' $bRowCount,1
for ($i = 0; $i -lt $bRowCount; $i++) {
    $bArr[$i,0] = $bValues[$i]
}
$ws.Range("B18:B89").Value2 = $bArr
